# Update 1.5.7: Bug fixes and QoL
# - Rename "Импорт" sheet to "СТВТ"
# - Add a "Бригады" sheet between "СТВТ" and "Супервайзеры"
# - Rework the header row of "СТВТ": swap/retarget the "Класс напряжения*" /
#   "ТТ Коэффицент*" columns, turn the old "Супервайзер*" header into a plain
#   "Супервайзер" column and append a new "Бригада" column
# - Add data validation lists for the new "Класс напряжения" and "Бригада" columns
# - Clear the sample rows out of "Супервайзеры" (keep just the header)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the first sheet "Импорт" -> "СТВТ"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "СТВТ"

# ---------------------------------------------------------------------------
# 2. Rework the header row on "СТВТ"
#    old: Наименование* | Статус* | Супервайзер*     | Класс напряжения* | ТТ Коэффицент*
#    new: Наименование* | Статус* | Класс напряжения* | ТТ Коэффицент*    | Супервайзер | Бригада
# ---------------------------------------------------------------------------
$ws1.Range("C1").Value = "Класс напряжения*"
$ws1.Range("D1").Value = "ТТ Коэффицент*"
$ws1.Range("E1").Value = "Супервайзер"

$ws1.Range("F1").Value = "Бригада"
$ws1.Range("F1").Font.Bold = $true
$ws1.Range("F1").Font.Size = 12
$ws1.Columns("F:F").ColumnWidth = 15.33

# ---------------------------------------------------------------------------
# 3. Data validation: "Класс напряжения" column gets a fixed list
# ---------------------------------------------------------------------------
$ws1.Range("C2:C1048576").Validation.Add(3, 1, 1, '"6кВ, 10кВ, 20кВ, 35кВ"')

# Move/retarget the supervisor lookup validation from column C to column E
$ws1.Range("E2:E1048576").Validation.Add(3, 1, 1, "=Супервайзеры!`$A`$2:`$A`$1048576")

# New "Бригада" column looks up the new "Бригады" sheet
$ws1.Range("F2:F1048576").Validation.Add(3, 1, 1, "=Бригады!`$A`$2:`$A`$58")

# Selection on "СТВТ" moves to E2
$ws1.Range("E2").Select()

# ---------------------------------------------------------------------------
# 4. Insert the new "Бригады" sheet right after "СТВТ"
# ---------------------------------------------------------------------------
$brigady = $wb.Worksheets.Add($null, $ws1)
$brigady.Name = "Бригады"
$brigady.Range("A1").Value = "Бригады"
$brigady.Range("A1").Font.Bold = $true
$brigady.Activate()
$brigady.Range("A2").Select()

# ---------------------------------------------------------------------------
# 5. Clear the sample rows out of "Супервайзеры", keep only the header
# ---------------------------------------------------------------------------
$wsSuper = $wb.Worksheets.Item("Супервайзеры")
$wsSuper.Range("A2:A4").EntireRow.Delete()
